# Remove the <sup>N</sup> suffixes from "ImageContrastModifier" labels in
# column B (Algorithm) of the Benchmarks sheet, so that
# "ImageContrastModifier<sup>1</sup>", "ImageContrastModifier<sup>2</sup>" and
# "ImageContrastModifier<sup>3</sup>" all collapse into a single
# "ImageContrastModifier" string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -like "ImageContrastModifier*") {
        $cell.Value2 = "ImageContrastModifier"
    }
}

# Reflect the cell selection recorded in the saved file.
$ws.Range("C21").Select() | Out-Null
